$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 1 de Abril de 2020 a las 20:50'
$ws.Range("B4").Value = 207157
$ws.Range("C4").Value = 18627
$ws.Range("E4").Value = 193789
$ws.Range("G4").Value = 553
$ws.Range("H4").Value = 4606
$ws.Range("B27").Value = 3107
$ws.Range("C27").Value = 247
$ws.Range("E27").Value = 2109
$ws.Range("E32").Value = 2116
$ws.Range("G32").Value = 10
$ws.Range("H32").Value = 92
$ws.Range("A50").Value = 'Colombia'
$ws.Range("B50").Value = 1065
$ws.Range("C50").Value = 159
$ws.Range("D50").Value = 39
$ws.Range("E50").Value = 1009
$ws.Range("F50").Value = 35
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 17
$ws.Range("A51").Value = 'Serbia'
$ws.Range("B51").Value = 1060
$ws.Range("C51").Value = 160
$ws.Range("D51").Value = 42
$ws.Range("E51").Value = 990
$ws.Range("F51").Value = 62
$ws.Range("G51").Value = 5
$ws.Range("A52").Value = 'Argentina'
$ws.Range("B52").Value = 1054
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 248
$ws.Range("E52").Value = 778
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 28
$ws.Range("A53").Value = 'Singapur'
$ws.Range("B53").Value = 1000
$ws.Range("C53").Value = 74
$ws.Range("D53").Value = 245
$ws.Range("E53").Value = 752
$ws.Range("F53").Value = 24
$ws.Range("H53").Value = 3
$ws.Range("A54").Value = 'Croacia'
$ws.Range("B54").Value = 963
$ws.Range("C54").Value = 96
$ws.Range("D54").Value = 73
$ws.Range("E54").Value = 884
$ws.Range("F54").Value = 34
$ws.Range("H54").Value = 6
$ws.Range("B57").Value = 835
$ws.Range("C57").Value = 54
$ws.Range("D57").Value = 71
$ws.Range("E57").Value = 762
$ws.Range("A58").Value = 'Emiratos Arabes Unidos'
$ws.Range("B58").Value = 814
$ws.Range("C58").Value = 150
$ws.Range("D58").Value = 61
$ws.Range("E58").Value = 745
$ws.Range("F58").Value = 2
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 8
$ws.Range("A59").Value = 'Estonia'
$ws.Range("B59").Value = 779
$ws.Range("C59").Value = 34
$ws.Range("D59").Value = 33
$ws.Range("E59").Value = 741
$ws.Range("F59").Value = 15
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 5
$ws.Range("A60").Value = 'Hong Kong'
$ws.Range("B60").Value = 765
$ws.Range("C60").Value = 50
$ws.Range("D60").Value = 147
$ws.Range("E60").Value = 614
$ws.Range("F60").Value = 5
$ws.Range("H60").Value = 4
$ws.Range("A61").Value = 'Crucero'
$ws.Range("B61").Value = 712
$ws.Range("D61").Value = 603
$ws.Range("E61").Value = 98
$ws.Range("F61").Value = 15
$ws.Range("H61").Value = 11
$ws.Range("A62").Value = 'Egipto'
$ws.Range("B62").Value = 710
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 157
$ws.Range("E62").Value = 507
$ws.Range("F62").Value = 0
$ws.Range("H62").Value = 46
$ws.Range("A63").Value = 'Nueva Zelanda'
$ws.Range("B63").Value = 708
$ws.Range("C63").Value = 61
$ws.Range("D63").Value = 83
$ws.Range("E63").Value = 624
$ws.Range("F63").Value = 2
$ws.Range("H63").Value = 1
$ws.Range("A64").Value = 'Irak'
$ws.Range("B64").Value = 694
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 170
$ws.Range("E64").Value = 474
$ws.Range("H64").Value = 50
$ws.Range("A65").Value = 'Ucrania'
$ws.Range("B65").Value = 669
$ws.Range("C65").Value = 24
$ws.Range("D65").Value = 10
$ws.Range("E65").Value = 642
$ws.Range("F65").Value = 0
$ws.Range("H65").Value = 17
$ws.Range("F86").Value = 14
$ws.Range("A91").Value = 'Afganistan'
$ws.Range("B91").Value = 239
$ws.Range("C91").Value = 65
$ws.Range("D91").Value = 5
$ws.Range("E91").Value = 230
$ws.Range("F91").Value = 0
$ws.Range("H91").Value = 4
$ws.Range("A92").Value = 'San Marino'
$ws.Range("B92").Value = 236
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 13
$ws.Range("E92").Value = 197
$ws.Range("F92").Value = 16
$ws.Range("H92").Value = 26
$ws.Range("A93").Value = 'Camerun'
$ws.Range("B93").Value = 233
$ws.Range("C93").Value = 40
$ws.Range("D93").Value = 10
$ws.Range("E93").Value = 217
$ws.Range("F93").Value = 0
$ws.Range("H93").Value = 6
$ws.Range("A94").Value = 'Vietnam'
$ws.Range("B94").Value = 218
$ws.Range("C94").Value = 6
$ws.Range("D94").Value = 63
$ws.Range("E94").Value = 155
$ws.Range("H94").Value = 0
$ws.Range("A95").Value = 'Cuba'
$ws.Range("B95").Value = 212
$ws.Range("C95").Value = 26
$ws.Range("D95").Value = 12
$ws.Range("E95").Value = 194
$ws.Range("H95").Value = 6
$ws.Range("A96").Value = 'Oman'
$ws.Range("B96").Value = 210
$ws.Range("C96").Value = 18
$ws.Range("D96").Value = 34
$ws.Range("E96").Value = 175
$ws.Range("F96").Value = 3
$ws.Range("H96").Value = 1
$ws.Range("E107").Value = 122
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 3
$ws.Range("B126").Value = 68
$ws.Range("C126").Value = 8
$ws.Range("E126").Value = 67
